$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated coin market data as text values (preserving exact formatting)
$updates = @{
    "D2" = "303.60"
    "E2" = "3.07%"
    "D3" = "33.59"
    "E3" = "8.17%"
    "D4" = "5.156"
    "E4" = "4.53%"
    "D5" = "0.07820"
    "E5" = "6.15%"
    "D6" = "2.399"
    "E6" = "4.56%"
    "D7" = "8.031"
    "E7" = "4.37%"
    "B8" = "MXToken"
    "C8" = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
    "D8" = "0.9338"
    "E8" = "2.28%"
    "B9" = "LiechtensteinCryptoassetsExchange"
    "C9" = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
    "D9" = "0.09816"
    "E9" = "17.32%"
    "B10" = "WazirX"
    "C10" = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
    "D10" = "0.1787"
    "E10" = "5.86%"
    "B11" = "MandalaExchangeToken"
    "C11" = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
    "D11" = "0.08509"
    "E11" = "3.23%"
    "B12" = "BitrueCoin"
    "C12" = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
    "D12" = "0.03349"
    "E12" = "7.38%"
    "B13" = "BitMartToken"
    "C13" = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
    "D13" = "0.09926"
    "E13" = "-1.50%"
    "B14" = "BitForexToken"
    "C14" = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
    "D14" = "0.001487"
    "E14" = "-1.44%"
    "B15" = "TigerCash"
    "C15" = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
    "D15" = "0.005801"
    "E15" = "0.43%"
    "B16" = "LEO"
    "C16" = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
    "D16" = "3.466"
    "E16" = "-0.49%"
    "B17" = "GateToken"
    "C17" = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
    "D17" = "3.914"
    "E17" = "4.16%"
    "D18" = "2.161"
    "E18" = "3.93%"
    "E19" = "1.15%"
    "D20" = "0.1342"
    "E20" = "3.00%"
    "D21" = "4.282"
    "E21" = "7.82%"
    "D22" = "0.2294"
    "E22" = "9.28%"
    "D23" = "0.04647"
    "E23" = "2.20%"
    "D24" = "0.001222"
    "E24" = "1.02%"
    "D25" = "0.004408"
    "E25" = "1.65%"
    "D26" = "0.0001295"
    "E26" = "-0.45%"
    "D27" = "0.0003391"
    "E27" = "-0.08%"
    "D39" = "0.01740"
    "E39" = "8.13%"
    "D40" = "0.04821"
    "E40" = "8.40%"
    "D41" = "0.007738"
    "E41" = "5.54%"
    "D42" = "0.009790"
    "E42" = "10.89%"
    "D43" = "0.1411"
    "E43" = "6.47%"
    "D44" = "0.002093"
    "E44" = "1.55%"
    "D45" = "0.009126"
    "E45" = "-0.69%"
    "D46" = "0.00006099"
    "E46" = "1.21%"
    "D47" = "0.00000000749"
    "E47" = "-0.15%"
    "D48" = "2.794"
    "E48" = "24.67%"
    "D49" = "0.001998"
    "E49" = "-31.08%"
    "D50" = "0.00002098"
    "E50" = "-0.15%"
    "D51" = "0.0001998"
    "E51" = "-0.15%"
}

foreach ($cell in $updates.Keys) {
    $range = $ws.Range($cell)
    $range.NumberFormat = "@"
    $range.Value = $updates[$cell]
}
